# Update cryptocurrency price/volume data per upstream GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking text (e.g. "0.999", "11.10").
# Left as the default "General" format, Excel would silently reinterpret
# these as real numbers (dropping trailing zeros, using scientific notation,
# introducing floating-point noise). Force a Text format first so the
# literal string from the source feed is preserved exactly, matching the
# original inline-string cells.
$textCells = @("D4", "D5", "D6", "D9", "D11", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D38", "D42", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.911.18"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "3.836.09"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "600.08"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "167.90"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").Value = "3.835.03"
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("D15").Value = "4.474.08"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "3.826.87"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").Value = "67.913.31"
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("D18").Value = "18.32"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "11.10"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "465.81"
$ws.Range("E22").Value = "  -5.97%  "
$ws.Range("D23").Value = "0.733"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "0.0000161"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").Value = "82.89"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "12.11"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "2.96"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "3.980.49"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "31.40"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "2.31"
$ws.Range("E34").Value = "  -6.33%  "
$ws.Range("D35").Value = "9.55"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "3.794.73"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").Value = "3.66"
$ws.Range("E38").Value = "  +10.30%  "
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -5.15%  "
$ws.Range("E44").Value = "  -7.01%  "
$ws.Range("D45").Value = "420.31"
$ws.Range("E45").Value = "  -4.46%  "
$ws.Range("D46").Value = "8.73"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "0.000294"
$ws.Range("E47").Value = "  +5.72%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "46.90"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "26.42"
$ws.Range("E50").Value = "  +4.51%  "
$ws.Range("D51").Value = "142.26"
$ws.Range("E51").Value = "  -0.84%  "
